$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# A new trade row is being inserted above the current row 5, pushing the
# existing rows 5-7 down to rows 6-8. Shift the data downward first
# (bottom-up, to avoid clobbering source rows before they're copied),
# preserving each cell's formatting exactly via per-cell Copy.
$cols = 1,2,3,4,5,6,7,9
foreach ($r in 7,6,5) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Copy($ws.Cells.Item($r + 1, $c))
    }
    $ws.Cells.Item($r + 1, 10).Formula = "=Index!`$C`$2"
}

# Write the new trade into row 5 (formatting matches the other data rows:
# only column A carries the date number format).
$ws.Cells.Item(5, 1).Value = 46062
$ws.Cells.Item(5, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(5, 2).Value = "NSE"
$ws.Cells.Item(5, 3).Value = "Buy"
$ws.Cells.Item(5, 4).Value = 100
$ws.Cells.Item(5, 5).Value = 37.8099
$ws.Cells.Item(5, 6).Value = 3799.99
$ws.Cells.Item(5, 7).Value = "CN#252611665409"
$ws.Cells.Item(5, 9).Value = 19
$ws.Cells.Item(5, 10).Formula = "=Index!`$C`$2"
